$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D13").Value = 0.4
$ws.Range("D14").Value = 0.6
$ws.Range("D15").Value = 0.8
$ws.Range("D16").Value = 0.8
$ws.Range("D17").Value = 0.8
$ws.Range("D18").Value = 0.8
$ws.Range("D19").Value = 0.8

$ws.Range("C22").Select()
